# "add more pearson r calcs"
# Populate Sheet5 with the "ratio" (E) and "earning" (F) columns, copied
# from the England rows of Sheet4 (ratio table) and Sheet3 (earnings
# table) respectively, and restore the per-sheet selections / active tab
# that resulted from that work session.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws5 = $wb.Worksheets.Item("Sheet5")

# ---------------------------------------------------------------------
# New column headers for Sheet5
# ---------------------------------------------------------------------
$ws5.Range("E1").Value = "ratio"
$ws5.Range("F1").Value = "earning"

# ---------------------------------------------------------------------
# Copy the cell formatting for the two new data columns from the source
# tables (England row) so the styles match exactly what was pasted in
# by the original author: Sheet4!C3 (ratio, 2002-2019) / Sheet4!U3
# (ratio, 2020) and Sheet3!C3 (earnings, 2002-2018) / Sheet3!T3
# (earnings, 2019-2020).
# ---------------------------------------------------------------------
$ws4.Range("C3").Copy()
$ws5.Range("E2:E19").PasteSpecial(-4122)

$ws4.Range("U3").Copy()
$ws5.Range("E20").PasteSpecial(-4122)

$ws3.Range("C3").Copy()
$ws5.Range("F2:F18").PasteSpecial(-4122)

$ws3.Range("T3").Copy()
$ws5.Range("F19:F20").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Values: "ratio" column (E) = lower quartile house price / earnings
# ratio for England, "earning" column (F) = lower quartile earnings
# for England.
# ---------------------------------------------------------------------
$ratio = @{2=4.41;3=5.17;4=6.16;5=6.7;6=7.15;7=7.26;8=6.96;9=6.41;10=6.76;11=6.64;12=6.61;13=6.51;14=6.73;15=6.94;16=7.05;17=7.15;18=7.18;19=7.05;20=7.01}
$earning = @{2=14755;3=15293;4=15909;5=16427;6=16645;7=17227;8=17968;9=18395;10=18495;11=18528;12=18920;13=19215;14=19317;15=19583;16=20141;17=20569;18=21165;19=21985;20=22813}

foreach ($row in 2..20) {
    $ws5.Range("E$row").Value = $ratio[$row]
    $ws5.Range("F$row").Value = $earning[$row]
}

# ---------------------------------------------------------------------
# Selections left behind on each sheet from the editing session. Sheet5
# (where the new data was added) is activated last so it becomes the
# workbook's active tab.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("D23").Select()

$ws3.Activate()
$ws3.Range("C3:U3").Select()

$ws6 = $wb.Worksheets.Item("Sheet6")
$ws6.Activate()
$ws6.Range("C2:C24").Select()

$ws7 = $wb.Worksheets.Item("Sheet7")
$ws7.Activate()
$ws7.Range("M19").Select()

$ws5.Activate()
$ws5.Range("F2:F20").Select()
